$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Populate new rows 335-366 (values), in the same order the source
# workbook author entered them (this governs shared-string insertion order)
$ws.Cells.Item(336, 1).Value = 2
$ws.Cells.Item(336, 2).Value = "荪"
$ws.Cells.Item(337, 1).Value = 2
$ws.Cells.Item(337, 2).Value = "煮饭233"
$ws.Cells.Item(338, 1).Value = 2
$ws.Cells.Item(338, 2).Value = "千里未来"
$ws.Cells.Item(339, 1).Value = 2
$ws.Cells.Item(339, 2).Value = "haooo"
$ws.Cells.Item(340, 1).Value = 2
$ws.Cells.Item(340, 2).Value = "闵佳瑶"
$ws.Cells.Item(341, 1).Value = 2
$ws.Cells.Item(341, 2).Value = "阿婆朱"
$ws.Cells.Item(342, 1).Value = 2
$ws.Cells.Item(342, 2).Value = "油腻的Wink😉"
$ws.Cells.Item(343, 1).Value = 1
$ws.Cells.Item(343, 2).Value = "带你嫖娼的馆长"
$ws.Cells.Item(344, 1).Value = 2
$ws.Cells.Item(344, 2).Value = "Obento"
$ws.Cells.Item(345, 1).Value = 2
$ws.Cells.Item(345, 2).Value = "银子舟"
$ws.Cells.Item(346, 1).Value = 2
$ws.Cells.Item(346, 2).Value = "大户爱"
$ws.Cells.Item(347, 1).Value = 2
$ws.Cells.Item(347, 2).Value = "冥渊"
$ws.Cells.Item(348, 1).Value = 2
$ws.Cells.Item(348, 2).Value = "Tony"
$ws.Cells.Item(349, 1).Value = 2
$ws.Cells.Item(349, 2).Value = "蜜汁帅气小彩蛋"
$ws.Cells.Item(350, 1).Value = 2
$ws.Cells.Item(350, 2).Value = "abc"
$ws.Cells.Item(351, 1).Value = 2
$ws.Cells.Item(351, 2).Value = "陌云铃"
$ws.Cells.Item(352, 1).Value = 2
$ws.Cells.Item(352, 2).Value = "arT"
$ws.Cells.Item(353, 1).Value = 2
$ws.Cells.Item(353, 2).Value = "忆丶醉"
$ws.Cells.Item(354, 1).Value = 2
$ws.Cells.Item(354, 2).Value = "宇宙空间法师"
$ws.Cells.Item(355, 1).Value = 2
$ws.Cells.Item(355, 2).Value = "玛瑙星mnx"
$ws.Cells.Item(356, 1).Value = 2
$ws.Cells.Item(356, 2).Value = "幻海"
$ws.Cells.Item(357, 1).Value = 2
$ws.Cells.Item(357, 2).Value = "陌云铃"
$ws.Cells.Item(358, 1).Value = 2
$ws.Cells.Item(358, 2).Value = "Marvin"
$ws.Cells.Item(359, 1).Value = 2
$ws.Cells.Item(359, 2).Value = "神烦"
$ws.Cells.Item(360, 1).Value = 2
$ws.Cells.Item(360, 2).Value = "Joke"
$ws.Cells.Item(361, 1).Value = 2
$ws.Cells.Item(361, 2).Value = "有酸萝卜别吃"
$ws.Cells.Item(362, 1).Value = 2
$ws.Cells.Item(362, 2).Value = "m子会梦见花之暴君吗"
$ws.Cells.Item(363, 1).Value = 1
$ws.Cells.Item(363, 2).Value = "刘云金"
$ws.Cells.Item(335, 1).Value = 2
$ws.Cells.Item(335, 2).Value = "天湮"
$ws.Cells.Item(343, 3).Value = "支持国产游戏，吾辈义不容辞！————读完本科自称硕士的沙雕的高质量网友"
$ws.Cells.Item(364, 1).Value = 2
$ws.Cells.Item(364, 2).Value = "蟹老板的老公"
$ws.Cells.Item(365, 1).Value = 2
$ws.Cells.Item(365, 2).Value = "小小书童"
$ws.Cells.Item(363, 3).Value = "我知道我们没有缘分，这也不是我想要的结果，事与愿违，曾经我们都想好好的，即使现在我也是想要跟你好好的，虽然你看不见，但还是想留下这么一句话（可能会有人看到会觉得挺可笑的 确实挺幼稚的）十年的恋情，不希望给你带来生活的压力， 希望你快快乐乐的 这也是我一直给你说的 也是一直希望的—GQQ"
$ws.Cells.Item(366, 1).Value = 2
$ws.Cells.Item(366, 2).Value = "雪华"
$ws.Cells.Item(366, 4).Value = "#00BFFF"

# --- Match formatting of the existing data rows (style index 3 = centered,
# bold-ish font used for B/C/D data cells further down the sheet)
$ws.Range("B334").Copy()
$ws.Range("B335:B366").PasteSpecial(-4122)
$ws.Range("C302").Copy()
$ws.Range("C343").PasteSpecial(-4122)
$ws.Range("C343").Copy()
$ws.Range("C363").PasteSpecial(-4122)
$ws.Range("D15").Copy()
$ws.Range("D366").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update the view: scroll target + new active selection (matches the
# commit's final cursor position after appending the new rows)
$ws.Range("G367").Select()
